$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Utiliza Facebook Sans, una tipografía Sans-serif clara y legible, priorizando la coherencia en sus interfaces. [3]"
$ws.Range("C2").Value = "Presenta una combinación de fuentes limpias para texto estándar, con opciones creativas para historias y gráficos. [5, 8]"
$ws.Range("D2").Value = "Muestra tipografías audaces y dinámicas, a menudo con efectos visuales, reflejando su naturaleza juvenil y enérgica. [2, 4, 14, 17]"
$ws.Range("E2").Value = "La tipografía es funcional pero carece de un estilo distintivo, usando fuentes estándar sin mucha personalidad ni jerarquía visual clara. [11, 20]"
$ws.Range("F2").Value = "Considerar implementar una tipografía más moderna y única para diferenciar la marca visualmente y mejorar la jerarquía."
$ws.Range("B3").Value = "Dominan los tonos azules clásicos y blanco, proyectando profesionalismo, confianza y una sensación de familiaridad. [15]"
$ws.Range("C3").Value = "Se caracteriza por degradados vibrantes y colores cálidos, evocando creatividad, alegría y una estética moderna y atractiva. [8]"
$ws.Range("D3").Value = "Emplea una paleta de colores oscuros y neón, generando una atmósfera moderna, juvenil y energéticamente estimulante. [1]"
$ws.Range("E3").Value = "La paleta de colores es básica y monocromática, lo que resulta en una apariencia algo genérica, plana y sin vida. [7, 11]"
$ws.Range("F3").Value = "Introducir una paleta de colores más rica y atractiva que refuerce la identidad de la plataforma y el contraste."
$ws.Range("B4").Value = "Adopta un tono predominantemente semi-formal, adecuado para conexiones personales, profesionales y comunicados oficiales. [13, 19]"
$ws.Range("C4").Value = "Mantiene un estilo informal y visualmente atractivo, fomentando la expresión personal y la interacción casual entre usuarios. [5, 8]"
$ws.Range("D4").Value = "Opera en un registro marcadamente informal y lúdico, celebrando la espontaneidad y el entretenimiento sin filtros. [2, 14]"
$ws.Range("E4").Value = "El lenguaje y el diseño son neutrales, cayendo en un rango formal que puede percibirse como poco cercano para el usuario. [6, 7]"
$ws.Range("F4").Value = "Evaluar un enfoque más informal en el lenguaje y diseño para fomentar mayor cercanía y participación activa de la audiencia."
$ws.Range("B5").Value = "Usa iconos planos y reconocibles universalmente, junto a un logo 'f' icónico y reacciones emotivas en el contenido. [12, 13, 15]"
$ws.Range("C5").Value = "Destaca por sus iconos minimalistas, el logo de la cámara fotográfica retro y elementos visuales para historias y reels. [5, 8]"
$ws.Range("D5").Value = "Presenta iconos dinámicos y un emblema vibrante que resuenan con su audiencia joven y la cultura de videos cortos. [1]"
$ws.Range("E5").Value = "Los iconos son genéricos y poco memorables, careciendo de un estilo unificado o elementos emblemáticos propios. [7]"
$ws.Range("F5").Value = "Desarrollar un conjunto de iconos personalizados y un emblema único que reflejen mejor la personalidad del sitio."
$ws.Range("B6").Value = "Ofrece múltiples funciones de accesibilidad, incluyendo lectores de pantalla y opciones de contraste para usuarios. [12]"
$ws.Range("C6").Value = "Ha mejorado su accesibilidad con texto alternativo para imágenes y subtítulos, aunque tiene áreas de oportunidad. [5, 8, 18]"
$ws.Range("D6").Value = "Continúa mejorando características de accesibilidad como subtítulos automáticos y ajustes de texto para inclusión plena. [1]"
$ws.Range("E6").Value = "La accesibilidad es básica, con limitaciones en opciones de contraste, texto alternativo y navegación por teclado, dificultando el acceso. [10, 11]"
$ws.Range("F6").Value = "Implementar un conjunto robusto de características de accesibilidad para garantizar un uso inclusivo para todos."
$ws.Range("B7").Value = "Posee una barra de navegación superior clara con accesos directos intuitivos a inicio, amigos y notificaciones. [19]"
$ws.Range("C7").Value = "Cuenta con una barra inferior prominente para inicio, búsqueda, reels, tienda y perfil, facilitando el acceso rápido. [5]"
$ws.Range("D7").Value = "Su navegación se centra en un feed vertical infinito, con botones clave de inicio, seguir, crear y perfil. [1]"
$ws.Range("E7").Value = "Los botones de navegación son estándar y su jerarquía no es siempre clara, lo que puede confundir al usuario. [6, 7, 10]"
$ws.Range("F7").Value = "Rediseñar la navegación para que los botones importantes sean más intuitivos, visibles y con una jerarquía clara."
$ws.Range("B8").Value = "Estructura el contenido en un feed central y barras laterales para amigos/grupos, manteniendo la información organizada. [19]"
$ws.Range("C8").Value = "Organiza el contenido de forma visual, priorizando cuadrículas de fotos, historias efímeras y reels de videos cortos. [5]"
$ws.Range("D8").Value = "Su organización se basa en un algoritmo que personaliza el feed, presentando videos relevantes sin estructura rígida. [1]"
$ws.Range("E8").Value = "La organización del contenido es lineal y predecible, careciendo de dinamismo o personalización para el usuario. [7, 11]"
$ws.Range("F8").Value = "Explorar nuevas formas de organizar el contenido, quizás con personalización, para mejorar la experiencia del usuario."
$ws.Range("B9").Value = "Ofrece Marketplace, Grupos, Eventos y Páginas, enriqueciendo la interacción y expandiendo las funcionalidades. [19]"
$ws.Range("C9").Value = "Incluye Stories, Reels, IGTV y la posibilidad de comprar directamente desde la aplicación, diversificando la experiencia. [5]"
$ws.Range("D9").Value = "Destaca por sus herramientas de edición de video, filtros creativos, efectos AR y duetos interactivos para usuarios. [1]"
$ws.Range("E9").Value = "Las características adicionales son mínimas o inexistentes, limitando la profundidad y la variedad de interacción del usuario."
$ws.Range("F9").Value = "Desarrollar e integrar características adicionales innovadoras que aporten valor y diferenciación a la plataforma."
$ws.Range("B10").Value = "Proporciona secciones de ayuda extensas y guías paso a paso para todas sus funcionalidades y herramientas complejas. [19]"
$ws.Range("C10").Value = "Ofrece tutoriales integrados para nuevas funciones y un centro de ayuda robusto con preguntas frecuentes claras. [5, 8]"
$ws.Range("D10").Value = "Brinda guías visuales y consejos en la aplicación, ayudando a los usuarios a dominar rápidamente sus herramientas creativas. [1]"
$ws.Range("E10").Value = "Las instrucciones son escasas y no siempre fáciles de encontrar, lo que dificulta el aprendizaje de nuevas funciones."
$ws.Range("F10").Value = "Crear tutoriales claros, accesibles y bien organizados, junto con un centro de ayuda intuitivo para usuarios."
$ws.Range("A11").Value = "Conclusion"
$ws.Range("B11").Value = "Una plataforma madura que equilibra comunicación personal con herramientas empresariales, manteniendo su relevancia global. [19]"
$ws.Range("C11").Value = "Se consolidó como líder visual, evolucionando constantemente para satisfacer las demandas de creadores y usuarios. [5]"
$ws.Range("D11").Value = "Revolucionó el contenido de video corto, creando una cultura vibrante y una experiencia de usuario altamente adictiva. [1]"
$ws.Range("E11").Value = "Este sitio ofrece una base sólida, pero con un potencial inmenso para mejorar su atractivo visual y funcionalidad. [6, 7]"
$ws.Range("F11").Value = "En general, el sitio podría beneficiarse enormemente de una renovación estética y funcional integral para destacar."
